## Aula_03_BD.pptx - "Modelos conceitual, logico e fisico finalizados"
##
## Slide 5 ("Um banco de dados e um conjunto de tabelas relacionadas..."),
## shape "Espaco Reservado para Conteudo 1": the sentence
##   "... tabela e um conjunto de informacoes sobre uma entidade ..."
## is reworded to
##   "... tabela e um conjunto de dados sobre uma entidade ..."
## i.e. "de informacoes " is replaced by "de dados ".

$p = $ppt.ActivePresentation
$slide = $p.Slides.Item(5)
$shape = $slide.Shapes.Item(1)
$tr = $shape.TextFrame.TextRange

# Locate the sentence unambiguously (this phrase only occurs once in the
# shape) so we don't have to hard-code character offsets.
$anchor = $tr.Find("conjunto de informa")

# "conjunto " (8 letters + trailing space) immediately precedes the chunk
# we need to retype; the next unique phrase "sobre uma entidade" marks
# where the replaced chunk ends.
$middleStart = $anchor.Start + 9
$tailAnchor = $tr.Find("sobre uma entidade")
$middleLen = $tailAnchor.Start - $middleStart

$toRetype = $tr.Characters($middleStart, $middleLen)
$toRetype.Text = "de dados "
